# Add the new output-path string (new shared string, uniqueCount 39 -> 40)
# and point the three "Output:" value cells (Message!A10,
# CypherOutput_Message!A10, StatOutput_Message!A20) at it, replacing the
# old shared-string reference (index 31) with the new one (index 39).

$wb = $excel.ActiveWorkbook

$newPath = "C:\Katalon_mastercopy520\DataCommons_Automation\OutputFiles\TC04_Canine_Filter_Breed-BassHnd_Neo4jData.xlsx"

$ws = $wb.Worksheets.Item("Message")
$ws.Range("A10").Value = $newPath

$ws = $wb.Worksheets.Item("CypherOutput_Message")
$ws.Range("A10").Value = $newPath

$ws = $wb.Worksheets.Item("StatOutput_Message")
$ws.Range("A20").Value = $newPath
